$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "RelivePos" values (186,6.89,88 -> 186,0,88) in column E (rows 2-4)
$ws.Range("E2").Value = "186,0,88"
$ws.Range("E3").Value = "186,0,88"
$ws.Range("E4").Value = "186,0,88"

# Update the active selection to F7
$ws.Range("F7").Select()
